# Weekly update: insert a new price record at row 10 for
# "Terminal Hortofrutícola Agro Chillán - Berenjena", pushing the
# existing rows 10-37 down to 11-38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 10 (shifts old rows 10..37 down to 11..38)
$ws.Rows.Item(10).Insert()

# Populate the new row 10 with the latest market data
$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(10, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(10, 3).Value = "Ñuble"
$ws.Cells.Item(10, 4).Value = 44764
$ws.Cells.Item(10, 5).Value = 16
$ws.Cells.Item(10, 6).Value = 100112001
$ws.Cells.Item(10, 7).Value = "Berenjena"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 60
$ws.Cells.Item(10, 11).Value = 12000
$ws.Cells.Item(10, 12).Value = 13000
$ws.Cells.Item(10, 13).Value = 12500
$ws.Cells.Item(10, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(10, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(10, 16).Value = 208
$ws.Cells.Item(10, 17).Value = 60
$ws.Cells.Item(10, 18).Value = "Hortaliza"

# Match the date cell number format used by the other rows in column D
$ws.Cells.Item(10, 4).NumberFormat = $ws.Cells.Item(11, 4).NumberFormat
